$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Ajo" (Chino / Primera) at the
# Terminal Hortofrutícola Agro Chillán market. Insert a fresh row at 117,
# pushing every existing row from 117 down to 247.
$ws.Rows("117").Insert()

$ws.Range("A117").Value = 7
$ws.Range("B117").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C117").Value = "Ñuble"
$ws.Range("D117").Value = 44789
$ws.Range("E117").Value = 16
$ws.Range("F117").Value = 100112003
$ws.Range("G117").Value = "Ajo"
$ws.Range("H117").Value = "Chino"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 100
$ws.Range("K117").Value = 24000
$ws.Range("L117").Value = 25000
$ws.Range("M117").Value = 24500
$ws.Range("N117").Value = "$/caja 10 kilos"
$ws.Range("O117").Value = "China"
$ws.Range("P117").Value = 2450
$ws.Range("Q117").Value = 10
$ws.Range("R117").Value = "Hortaliza"
